{"js": "const replacements = [\n  [\"2024-06-25 Tuesday\", \"2024-06-26 Wednesday\"],\n  [\"925\u00d79=8325\", \"835\u00d77=5845\"],\n  [\"287\u00d74=1148\", \"995\u00d79=8955\"],\n  [\"992\u00d75=4960\", \"962\u00d76=5772\"],\n  [\"212\u00d77=1484\", \"764\u00d72=1528\"],\n  [\"930\u00d73=2790\", \"998\u00d79=8982\"],\n  [\"846\u00d79=7614\", \"325\u00d76=1950\"],\n  [\"260\u00d73=780\", \"669\u00d79=6021\"],\n  [\"165\u00d73=495\", \"599\u00d76=3594\"],\n  [\"325\u00d74=1300\", \"610\u00d79=5490\"],\n  [\"821\u00d74=3284\", \"443\u00d74=1772\"],\n  [\"413\u00d72=826\", \"699\u00d77=4893\"],\n  [\"935\u00d73=2805\", \"395\u00d79=3555\"],\n  [\"133\u00d76=798\", \"274\u00d74=1096\"],\n  [\"696\u00d73=2088\", \"648\u00d79=5832\"],\n  [\"621\u00d72=1242\", \"725\u00d73=2175\"],\n  [\"296\u00d76=1776\", \"122\u00d74=488\"],\n  [\"766\u00d74=3064\", \"849\u00d79=7641\"],\n  [\"832\u00d75=4160\", \"300\u00d73=900\"],\n  [\"782\u00d76=4692\", \"332\u00d74=1328\"],\n  [\"940\u00d77=6580\", \"870\u00d78=6960\"],\n  [\"312\u00d78=2496\", \"674\u00d72=1348\"],\n  [\"984\u00d79=8856\", \"822\u00d77=5754\"],\n  [\"776\u00d72=1552\", \"283\u00d77=1981\"],\n  [\"818\u00d79=7362\", \"916\u00d76=5496\"],\n  [\"798\u00d74=3192\", \"310\u00d72=620\"],\n];\n\nconst body = context.document.body;\n\n// Replace each old value with its new value. Every search string is\n// unique within the document (date line + one cell per multiplication\n// problem), so each search should resolve to exactly one hit.\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-Text '2024-06-25 Tuesday' '2024-06-26 Wednesday'\nReplace-Text '925\u00d79=8325' '835\u00d77=5845'\nReplace-Text '287\u00d74=1148' '995\u00d79=8955'\nReplace-Text '992\u00d75=4960' '962\u00d76=5772'\nReplace-Text '212\u00d77=1484' '764\u00d72=1528'\nReplace-Text '930\u00d73=2790' '998\u00d79=8982'\nReplace-Text '846\u00d79=7614' '325\u00d76=1950'\nReplace-Text '260\u00d73=780' '669\u00d79=6021'\nReplace-Text '165\u00d73=495' '599\u00d76=3594'\nReplace-Text '325\u00d74=1300' '610\u00d79=5490'\nReplace-Text '821\u00d74=3284' '443\u00d74=1772'\nReplace-Text '413\u00d72=826' '699\u00d77=4893'\nReplace-Text '935\u00d73=2805' '395\u00d79=3555'\nReplace-Text '133\u00d76=798' '274\u00d74=1096'\nReplace-Text '696\u00d73=2088' '648\u00d79=5832'\nReplace-Text '621\u00d72=1242' '725\u00d73=2175'\nReplace-Text '296\u00d76=1776' '122\u00d74=488'\nReplace-Text '766\u00d74=3064' '849\u00d79=7641'\nReplace-Text '832\u00d75=4160' '300\u00d73=900'\nReplace-Text '782\u00d76=4692' '332\u00d74=1328'\nReplace-Text '940\u00d77=6580' '870\u00d78=6960'\nReplace-Text '312\u00d78=2496' '674\u00d72=1348'\nReplace-Text '984\u00d79=8856' '822\u00d77=5754'\nReplace-Text '776\u00d72=1552' '283\u00d77=1981'\nReplace-Text '818\u00d79=7362' '916\u00d76=5496'\nReplace-Text '798\u00d74=3192' '310\u00d72=620'\n"}
